# Generate Report for Handoff
# Rotate the handed-off file's GUID/commit-hash identifiers and bump the
# recorded timestamps, across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "a23268c8-93ce-4e31-82af-342e58d1d30f"
$newGuid = "f11706e4-eeb2-43d4-bd78-9291a933d500"

$oldHash = "330c0a7d3068577afd852a26975dd98ca29e1be9"
$newHash = "096b46e925a0e380cce2081706862ffc8166dda7"

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newGuid + ".md"
$wsOverview.Range("B2").Value = "e2e\" + $newGuid + ".md"
$wsOverview.Range("G2").Value = "2016-08-15 12:52:53"
foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = "e2e\" + $newGuid + ".md"
}

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = $newGuid + ".md"
$wsZhCn.Range("G2").Value = $newGuid + "." + $newHash + ".zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-15 12:52:48"
foreach ($hl in $wsZhCn.Hyperlinks) {
    $hl.TextToDisplay = $newGuid + ".md"
}

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = $newGuid + ".md"
$wsDeDe.Range("G2").Value = $newGuid + "." + $newHash + ".de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-15 12:52:53"
foreach ($hl in $wsDeDe.Hyperlinks) {
    $hl.TextToDisplay = $newGuid + ".md"
}

Write-Host "Report regenerated for handoff."
